# Update 'want to go' counts (column F) across sheets, per upstream data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 234
$ws.Range("F3").Value = 1427
$ws.Range("F4").Value = 20083
$ws.Range("F6").Value = 314
$ws.Range("F7").Value = 1098
$ws.Range("F9").Value = 7576
$ws.Range("F10").Value = 515
$ws.Range("F12").Value = 267
$ws.Range("F13").Value = 38
$ws.Range("F15").Value = 119
$ws.Range("F17").Value = 235
$ws.Range("F19").Value = 1337
$ws.Range("F20").Value = 421
$ws.Range("F23").Value = 48
$ws.Range("F24").Value = 70
$ws.Range("F26").Value = 320
$ws.Range("F27").Value = 1106
$ws.Range("F29").Value = 19
$ws.Range("F30").Value = 183
$ws.Range("F31").Value = 5220
$ws.Range("F32").Value = 562
$ws.Range("F33").Value = 66
$ws.Range("F34").Value = 2848
$ws.Range("F36").Value = 87
$ws.Range("F38").Value = 12643
$ws.Range("F39").Value = 1336
$ws.Range("F40").Value = 82
$ws.Range("F41").Value = 26
$ws.Range("F42").Value = 55
$ws.Range("F43").Value = 264
$ws.Range("F44").Value = 365
$ws.Range("F45").Value = 4001

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 175

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 234
$ws.Range("F3").Value = 1427
$ws.Range("F4").Value = 20083
$ws.Range("F6").Value = 314
$ws.Range("F7").Value = 1098
$ws.Range("F9").Value = 7576
$ws.Range("F10").Value = 515
$ws.Range("F12").Value = 267
$ws.Range("F13").Value = 38
$ws.Range("F15").Value = 119
$ws.Range("F17").Value = 235
$ws.Range("F19").Value = 1337
$ws.Range("F20").Value = 421
$ws.Range("F23").Value = 48
$ws.Range("F24").Value = 70
$ws.Range("F26").Value = 320
$ws.Range("F27").Value = 1106
$ws.Range("F29").Value = 19
$ws.Range("F30").Value = 183
$ws.Range("F31").Value = 175
$ws.Range("F32").Value = 562
$ws.Range("F34").Value = 66
$ws.Range("F36").Value = 2848
$ws.Range("F38").Value = 87
$ws.Range("F40").Value = 12643
$ws.Range("F41").Value = 1336
$ws.Range("F42").Value = 82
$ws.Range("F43").Value = 26
$ws.Range("F44").Value = 55
$ws.Range("F45").Value = 264
$ws.Range("F46").Value = 365
$ws.Range("F47").Value = 4001
